# Append a new price-check column (H) to the LDLC price-tracking sheet,
# mirroring the most recent snapshot (column G) - "Update LDLC prices history".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell (H1): new timestamp, same header style as the other
#     timestamp columns (bold / bordered / centered = style of G1). ---
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "2026-01-27 21:09:45"

# --- Data rows 2-100: column H duplicates column G's numeric price for
#     this snapshot (read as Value2 to keep full numeric precision, then
#     bulk-assign the whole column in one shot). ---
$prices = $ws.Range("G2:G100").Value2()
$ws.Range("H2:H100").Value2 = $prices

# --- Trailer rows 101-204: column G holds an (empty) text cell for these
#     placeholder/out-of-stock rows; replicate that as an empty text cell
#     in column H too (leading "'" forces text type with an empty value,
#     then resetting the style keeps it unformatted just like column G). ---
$blankRange = $ws.Range("H101:H204")
$blankRange.Value = "'"
$blankRange.Style = "Normal"

Write-Output "done"
